# Add raw materials for steel production costs to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Match styling used by the other parameter label cells in column A
# (wrap text font) for the three new label cells, same as the row above.
$ws.Range("A25").Copy()
$ws.Range("A27:A29").PasteSpecial(-4122)

# New parameter rows: Lime, Carbon, Iron Ore Pellets (with per-site costs)
$ws.Range("A27").Value = "Lime (`$/metric tonne)"
$ws.Range("B27").Value = 155.34
$ws.Range("C27").Value = 141.51
$ws.Range("D27").Value = 155.34
$ws.Range("E27").Value = 169.18

$ws.Range("A28").Value = "Carbon (`$/metric tonne)"
$ws.Range("B28").Value = 218.74
$ws.Range("C28").Value = 264.94
$ws.Range("D28").Value = 229.74
$ws.Range("E28").Value = 245.14

$ws.Range("A29").Value = "Iron Ore Pellets (`$/metric tonne)"
$ws.Range("B29").Value = 230.52
$ws.Range("C29").Value = 292.2
$ws.Range("D29").Value = 239.32
$ws.Range("E29").Value = 270.2

# The "Iron Ore Pellets" label wraps onto two lines, so the row is taller.
$ws.Rows.Item(29).RowHeight = 32

# Update the active selection to match the author's final cursor position.
$ws.Range("H18").Select()
